$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The scraped source refreshes the crypto price/volume table. Column D ("Price")
# holds values as text (some look like plain numbers, e.g. "591.74", others use
# dots as thousands separators, e.g. "66.598.62"). Column E ("Volume(1h)") always
# holds a padded percentage string, e.g. "  +1.48%  ", which is never ambiguous.
#
# Assigning a plain numeric-looking string via .Value would make Excel coerce it
# into a real number, changing the cell's stored type. To keep these as text -
# matching the source data - prefix numeric-looking values with an apostrophe,
# exactly like typing '591.74 into a cell in the Excel UI.

function Set-TextValue($worksheet, $row, $col, $value) {
    if ($value -match '^-?\d+(\.\d+)?$') {
        $worksheet.Cells.Item($row, $col).Value = "'" + $value
    } else {
        $worksheet.Cells.Item($row, $col).Value = $value
    }
}

$changes = @(
    @{ Row=2; D='66.598.62'; E='  +1.48%  ' },
    @{ Row=3; D='3.497.09'; E='  +0.93%  ' },
    @{ Row=5; D='591.74'; E='  +1.72%  ' },
    @{ Row=6; D='168.64'; E='  +0.39%  ' },
    @{ Row=7; E='  +0.01%  ' },
    @{ Row=8; E='  +4.52%  ' },
    @{ Row=9; E='  +4.68%  ' },
    @{ Row=10; E='  +0.54%  ' },
    @{ Row=11; E='  +0.46%  ' },
    @{ Row=12; D='4.100.64'; E='  +1.10%  ' },
    @{ Row=13; E='  -0.19%  ' },
    @{ Row=14; E='  +2.33%  ' },
    @{ Row=15; E='  +1.55%  ' },
    @{ Row=16; D='66.631.11' },
    @{ Row=17; D='3.493.67'; E='  +0.78%  ' },
    @{ Row=18; E='  +1.21%  ' },
    @{ Row=19; E='  +2.00%  ' },
    @{ Row=20; D='393.60'; E='  +2.11%  ' },
    @{ Row=21; E='  -0.04%  ' },
    @{ Row=22; D='73.07'; E='  +2.16%  ' },
    @{ Row=23; E='  -0.12%  ' },
    @{ Row=24; D='0.533'; E='  +2.39%  ' },
    @{ Row=25; D='0.0000121'; E='  +1.16%  ' },
    @{ Row=26; D='10.16'; E='  +3.55%  ' },
    @{ Row=27; E='  -0.72%  ' },
    @{ Row=28; D='1.00'; E='  +0.45%  ' },
    @{ Row=29; E='  +1.76%  ' },
    @{ Row=30; E='  +0.47%  ' },
    @{ Row=31; E='  +1.36%  ' },
    @{ Row=32; D='23.80'; E='  +2.15%  ' },
    @{ Row=33; E='  +0.08%  ' },
    @{ Row=34; E='  +5.12%  ' },
    @{ Row=35; D='162.42'; E='  +1.40%  ' },
    @{ Row=36; E='  +0.43%  ' },
    @{ Row=37; E='  +2.82%  ' },
    @{ Row=38; D='6.79'; E='  +2.61%  ' },
    @{ Row=39; D='4.66'; E='  +4.45%  ' },
    @{ Row=40; E='  +1.73%  ' },
    @{ Row=41; D='0.0739'; E='  +0.67%  ' },
    @{ Row=42; D='26.81'; E='  -0.07%  ' },
    @{ Row=43; D='2.776.36'; E='  -0.79%  ' },
    @{ Row=44; D='42.92'; E='  -0.33%  ' },
    @{ Row=45; E='  +3.30%  ' },
    @{ Row=46; E='  +0.44%  ' },
    @{ Row=47; D='341.82'; E='  +1.38%  ' },
    @{ Row=48; E='  +1.13%  ' },
    @{ Row=49; D='34.03'; E='  +4.85%  ' },
    @{ Row=50; D='0.856'; E='  +2.85%  ' },
    @{ Row=51; E='  +1.67%  ' }
)

foreach ($change in $changes) {
    $r = $change.Row
    if ($change.ContainsKey('D')) {
        Set-TextValue $ws $r 4 $change.D
    }
    if ($change.ContainsKey('E')) {
        Set-TextValue $ws $r 5 $change.E
    }
}
